$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("profile")

# Insert a new column before column N (14th column), shifting existing
# columns N.. right by one.
$ws.Columns.Item(14).Insert()

# New header for the inserted column.
$ws.Cells.Item(1, 14).Value = "pro_usda_soil_order"

# Fill in the full USDA soil order name based on the (now shifted)
# abbreviated soil taxon column O.
for ($r = 2; $r -le 9; $r++) {
    $abbrev = $ws.Cells.Item($r, 15).Value()
    if ($abbrev -eq "Spo") {
        $ws.Cells.Item($r, 14).Value = "Spodosols"
    } elseif ($abbrev -eq "Ept") {
        $ws.Cells.Item($r, 14).Value = "Inceptisols"
    }
}

$ws.Activate()
$ws.Range("N10").Select()
